$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2925757.38
$ws.Range("C7").Value = -34.150215211611
$ws.Range("D7").Value = 2954
$ws.Range("E7").Value = 2954
$ws.Range("F7").Value = 990.4391943127962
$ws.Range("G7").Value = 5.573656316117237
